$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 and E2 text values
$ws.Range("D2").Value = "TC04_CDS_phs001524_SampleTumorStatus_NSD_TSVData.xlsx"
$ws.Range("E2").Value = "TC04_CDS_phs001524_SampleTumorStatus_NSD_WebData.xlsx"

# Update B3 query text - remove Tumor and Analyte Type columns
$ws.Range("B3").Value = "SELECT`n    DISTINCT (smp.sample_id) AS ""Sample ID"",`n    sp.participant_id AS ""Participant ID"", `n    s.study_name AS ""Study Name"",`n    s.phs_accession AS Accession`nFROM `n    df_participant sp`nJOIN `n    df_study s ON sp.""study.phs_accession"" = s.phs_accession`nJOIN `n    df_sample smp ON smp.""participant.study_participant_id"" = sp.study_participant_id`nJOIN`n    df_diagnosis d ON d.""participant.study_participant_id"" = sp.study_participant_id`nJOIN`n    df_program p ON p.program_acronym = s.""program.program_acronym""`nJOIN`n    df_file f1 ON f1.""sample.sample_id"" = smp.sample_id`nJOIN`n    df_genomic_info gi ON gi.""file.file_id"" = f1.file_id`nWHERE `n    s.phs_accession = 'phs001524' AND smp.sample_tumor_status = 'Not specified in data'`nORDER BY `n    smp.sample_id ASC`nLIMIT 100;"

# Clear D3, E3, D4, E4
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Update sheet view: scroll so row 3 is at top, select C3
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C3").Select()
